$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("D2").Value = 0.0001403349451720715
$ws.Range("E2").Value = 0.03105690097436309
$ws.Range("G2").Value = 0.002919693943113089
$ws.Range("H2").Value = 0.005493101663887501
$ws.Range("I2").Value = 0.008178039453923702
$ws.Range("J2").Value = 0.01082539837807417
$ws.Range("K2").Value = 0.001010153442621231
$ws.Range("D3").Value = 0.001920620910823345
$ws.Range("E3").Value = 0.0349030657671392
$ws.Range("G3").Value = 0.002867874689400196
$ws.Range("H3").Value = 0.007513747084885836
$ws.Range("I3").Value = 0.008498859126120806
$ws.Range("J3").Value = 0.01235072687268257
$ws.Range("K3").Value = 0.001014214474707842
$ws.Range("D4").Value = 0.001986748073250055
$ws.Range("E4").Value = 0.03560670325532556
$ws.Range("G4").Value = 0.002889716066420078
$ws.Range("H4").Value = 0.007650318555533886
$ws.Range("I4").Value = 0.00861393054947257
$ws.Range("J4").Value = 0.0124573796056211
$ws.Range("K4").Value = 0.001009044237434864
$ws.Range("D5").Value = 0.0002035470679402351
$ws.Range("E5").Value = 0.03206363087520003
$ws.Range("G5").Value = 0.002831997349858284
$ws.Range("H5").Value = 0.005671877413988113
$ws.Range("I5").Value = 0.008593371137976646
$ws.Range("J5").Value = 0.01137784495949745
$ws.Range("K5").Value = 0.001032428350299597
$ws.Range("D6").Value = 0.002898410893976688
$ws.Range("E6").Value = 0.1634307177737355
$ws.Range("G6").Value = 0.005984170362353325
$ws.Range("H6").Value = 0.01477164914831519
$ws.Range("I6").Value = 0.1164159486070275
$ws.Range("J6").Value = 0.0182621437124908
$ws.Range("K6").Value = 0.002305898815393448
$ws.Range("D8").Value = 0.0001403349451720715
$ws.Range("E8").Value = 0.03105690097436309
$ws.Range("G8").Value = 0.002919693943113089
$ws.Range("H8").Value = 0.005493101663887501
$ws.Range("I8").Value = 0.008178039453923702
$ws.Range("J8").Value = 0.01082539837807417
$ws.Range("K8").Value = 0.001010153442621231
$ws.Range("D9").Value = 0.001920620910823345
$ws.Range("E9").Value = 0.0349030657671392
$ws.Range("G9").Value = 0.002867874689400196
$ws.Range("H9").Value = 0.007513747084885836
$ws.Range("I9").Value = 0.008498859126120806
$ws.Range("J9").Value = 0.01235072687268257
$ws.Range("K9").Value = 0.001014214474707842
$ws.Range("D10").Value = 0.001986748073250055
$ws.Range("E10").Value = 0.03560670325532556
$ws.Range("G10").Value = 0.002889716066420078
$ws.Range("H10").Value = 0.007650318555533886
$ws.Range("I10").Value = 0.00861393054947257
$ws.Range("J10").Value = 0.0124573796056211
$ws.Range("K10").Value = 0.001009044237434864
$ws.Range("D11").Value = 0.0002035470679402351
$ws.Range("E11").Value = 0.03206363087520003
$ws.Range("G11").Value = 0.002831997349858284
$ws.Range("H11").Value = 0.005671877413988113
$ws.Range("I11").Value = 0.008593371137976646
$ws.Range("J11").Value = 0.01137784495949745
$ws.Range("K11").Value = 0.001032428350299597
$ws.Range("D12").Value = 0.002898410893976688
$ws.Range("E12").Value = 0.1634307177737355
$ws.Range("G12").Value = 0.005984170362353325
$ws.Range("H12").Value = 0.01477164914831519
$ws.Range("I12").Value = 0.1164159486070275
$ws.Range("J12").Value = 0.0182621437124908
$ws.Range("K12").Value = 0.002305898815393448
